# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1, styled like the other header cells (bold/bordered = same
# style as G1 "sum").
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Save values for rows 2-43: all 0 except rows 5 and 38 which are 1.
$saveValues = @{
    5  = 1
    38 = 1
}

for ($r = 2; $r -le 43; $r++) {
    if ($saveValues.ContainsKey($r)) {
        $ws.Cells.Item($r, 8).Value = $saveValues[$r]
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
